$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1278.6666
$ws.Range("J38").Value = 3363.5
$ws.Range("L38").Value = 10090.5
$ws.Range("N38").Value = -10834.5
$ws.Range("H39").Value = 618.25
$ws.Range("I39").Value = 124.333336
$ws.Range("K39").Value = 373.000008
$ws.Range("M39").Value = -77.00000799999998
$ws.Range("H40").Value = 9567.333000000001
$ws.Range("I40").Value = 9900
$ws.Range("J40").Value = 8902
$ws.Range("K40").Value = 9900
$ws.Range("L40").Value = 8902
$ws.Range("M40").Value = -9725
$ws.Range("N40").Value = -9252
$ws.Range("H41").Value = 235.6875
$ws.Range("I41").Value = 292.42856
$ws.Range("K41").Value = 292.42856
$ws.Range("M41").Value = 147.57144
$ws.Range("H43").Value = 3808.111
$ws.Range("I43").Value = 2738
$ws.Range("J43").Value = 4113.857
$ws.Range("K43").Value = 2738
$ws.Range("L43").Value = 4113.857
$ws.Range("M43").Value = -2669
$ws.Range("N43").Value = -4251.857
$ws.Range("H82").Value = 5028.8
$ws.Range("I82").Value = 1382
$ws.Range("J82").Value = 10499
$ws.Range("K82").Value = 4146
$ws.Range("L82").Value = 31497
$ws.Range("M82").Value = -3740
$ws.Range("N82").Value = -32309
$ws.Range("H85").Value = 5028.8
$ws.Range("I85").Value = 1382
$ws.Range("J85").Value = 10499
$ws.Range("K85").Value = 4146
$ws.Range("L85").Value = 31497
$ws.Range("M85").Value = -2742
$ws.Range("N85").Value = -34305
$ws.Range("H86").Value = 4666
$ws.Range("I86").Value = 5749
$ws.Range("K86").Value = 5749
$ws.Range("M86").Value = -4626
$ws.Range("H89").Value = 4666
$ws.Range("I89").Value = 5749
$ws.Range("K89").Value = 28745
$ws.Range("M89").Value = -23129
$ws.Range("H92").Value = 952.25
$ws.Range("I92").Value = 802.7143
$ws.Range("J92").Value = 1999
$ws.Range("K92").Value = 802.7143
$ws.Range("L92").Value = 1999
$ws.Range("M92").Value = 445.2857
$ws.Range("N92").Value = -4495
$ws.Range("H116").Value = 9550.629999999999
$ws.Range("I116").Value = 8603.130999999999
$ws.Range("J116").Value = 14998.75
$ws.Range("K116").Value = 8603.130999999999
$ws.Range("L116").Value = 14998.75
$ws.Range("M116").Value = -5161.130999999999
$ws.Range("N116").Value = -21882.75
$ws.Range("H135").Value = 2198.5
$ws.Range("I135").Value = 1110.75
$ws.Range("K135").Value = 9996.75
$ws.Range("M135").Value = -7461.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 50000000
$ws.Range("I37").Value = 50000000
$ws.Range("K37").Value = 50000000
$ws.Range("M37").Value = -49999727
$ws.Range("H61").Value = 3956.9412
$ws.Range("I61").Value = 2076.889
$ws.Range("K61").Value = 2076.889
$ws.Range("M61").Value = -1864.889
$ws.Range("H74").Value = 1654.4546
$ws.Range("I74").Value = 1619.9
$ws.Range("K74").Value = 1619.9
$ws.Range("M74").Value = -745.9000000000001
$ws.Range("H77").Value = 1654.4546
$ws.Range("I77").Value = 1619.9
$ws.Range("K77").Value = 8099.5
$ws.Range("M77").Value = -3731.5
$ws.Range("H136").Value = 3956.9412
$ws.Range("I136").Value = 2076.889
$ws.Range("K136").Value = 6230.667
$ws.Range("M136").Value = -3680.667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3419
$ws.Range("I31").Value = 1474
$ws.Range("K31").Value = 1474
$ws.Range("M31").Value = -1179
$ws.Range("H33").Value = 24146.166
$ws.Range("J33").Value = 30035
$ws.Range("L33").Value = 30035
$ws.Range("N33").Value = -30793
$ws.Range("H34").Value = 3419
$ws.Range("I34").Value = 1474
$ws.Range("K34").Value = 1474
$ws.Range("M34").Value = -1272
$ws.Range("H122").Value = 3437.6667
$ws.Range("J122").Value = 4019.2
$ws.Range("L122").Value = 12057.6
$ws.Range("N122").Value = -16957.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 685.2727
$ws.Range("I50").Value = 594.2222
$ws.Range("K50").Value = 1782.6666
$ws.Range("M50").Value = -1301.6666
$ws.Range("H53").Value = 685.2727
$ws.Range("I53").Value = 594.2222
$ws.Range("K53").Value = 1782.6666
$ws.Range("M53").Value = -1301.6666
$ws.Range("H76").Value = 12263
$ws.Range("I76").Value = 9735.6
$ws.Range("J76").Value = 24900
$ws.Range("K76").Value = 29206.8
$ws.Range("L76").Value = 74700
$ws.Range("M76").Value = -28823.8
$ws.Range("N76").Value = -75466
$ws.Range("H79").Value = 12263
$ws.Range("I79").Value = 9735.6
$ws.Range("J79").Value = 24900
$ws.Range("K79").Value = 29206.8
$ws.Range("L79").Value = 74700
$ws.Range("M79").Value = -27880.8
$ws.Range("N79").Value = -77352
$ws.Range("H122").Value = 160058.44
$ws.Range("J122").Value = 160058.44
$ws.Range("L122").Value = 1440525.96
$ws.Range("N122").Value = -1445425.96
$ws.Range("H131").Value = 37683044
$ws.Range("I131").Value = 59259750
$ws.Range("J131").Value = 23812300
$ws.Range("K131").Value = 177779250
$ws.Range("L131").Value = 71436900
$ws.Range("M131").Value = -177774210
$ws.Range("N131").Value = -71446980

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3954.4814
$ws.Range("I132").Value = 3335.6316
$ws.Range("K132").Value = 10006.8948
$ws.Range("M132").Value = -7476.8948

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3617.3076
$ws.Range("I16").Value = 2272.9167
$ws.Range("K16").Value = 2272.9167
$ws.Range("M16").Value = -2102.9167
$ws.Range("H93").Value = 1106.0769
$ws.Range("I93").Value = 339
$ws.Range("K93").Value = 339
$ws.Range("M93").Value = 909
$ws.Range("H132").Value = 2584.25
$ws.Range("I132").Value = 2443
$ws.Range("J132").Value = 2631.3333
$ws.Range("K132").Value = 7329
$ws.Range("L132").Value = 7893.999899999999
$ws.Range("M132").Value = -4799
$ws.Range("N132").Value = -12953.9999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 50000000
$ws.Range("I9").Value = 50000000
$ws.Range("K9").Value = 50000000
$ws.Range("M9").Value = -49999860
$ws.Range("H81").Value = 959.53845
$ws.Range("I81").Value = 959.53845
$ws.Range("K81").Value = 1919.0769
$ws.Range("M81").Value = -858.0769
$ws.Range("H84").Value = 959.53845
$ws.Range("I84").Value = 959.53845
$ws.Range("K84").Value = 9595.3845
$ws.Range("M84").Value = -4291.3845

